$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new worksheet "2022-Q4" right after "总计" (i.e. before
#    the current "2022-Q3" sheet), matching the workbook.xml sheet ordering
#    in the target diff.
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add($refSheet)
$q4.Name = "2022-Q4"

# Match the sheetPr/outlinePr defaults (summaryBelow/summaryRight) that
# every other sheet in the workbook carries (touching these properties is
# what makes the engine emit the otherwise-default <sheetPr> element).
$q4.Outline.SummaryRow = 1
$q4.Outline.SummaryColumn = -4152

# Header row (text labels)
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Columns B:G hold text values in the source data (fund codes with leading
# zeros, and numeric-looking strings that must stay text) - force text
# format before writing so Excel doesn't silently coerce them to numbers,
# then reset back to the Normal style so the cells don't carry a stray
# number-format style (matches the target, whose data cells are unstyled).
$q4.Range("B2:G6").NumberFormat = "@"

# Data rows
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "005802"
$q4.Range("C2").Value = "汇添富智能制造股票A"
$q4.Range("D2").Value = "24.50"
$q4.Range("E2").Value = "91.67"
$q4.Range("F2").Value = "3.77"
$q4.Range("G2").Value = "0.9236"
$q4.Range("H2").Value = 10

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "015797"
$q4.Range("C3").Value = "万家新能源主题混合C"
$q4.Range("D3").Value = "0.69"
$q4.Range("E3").Value = "93.18"
$q4.Range("F3").Value = "4.14"
$q4.Range("G3").Value = "0.0286"
$q4.Range("H3").Value = 6

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "015796"
$q4.Range("C4").Value = "万家新能源主题混合A"
$q4.Range("D4").Value = "0.31"
$q4.Range("E4").Value = "93.18"
$q4.Range("F4").Value = "4.14"
$q4.Range("G4").Value = "0.0128"
$q4.Range("H4").Value = 6

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "015196"
$q4.Range("C5").Value = "汇添富智能制造股票C"
$q4.Range("D5").Value = "0.02"
$q4.Range("E5").Value = "91.67"
$q4.Range("F5").Value = "3.77"
$q4.Range("G5").Value = "0.0008"
$q4.Range("H5").Value = 10

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "015197"
$q4.Range("C6").Value = "汇添富智能制造股票D"
$q4.Range("D6").Value = "0.01"
$q4.Range("E6").Value = "91.67"
$q4.Range("F6").Value = "3.77"
$q4.Range("G6").Value = "0.0004"
$q4.Range("H6").Value = 10

# Drop the temporary text number-format now that every text value is
# committed, so the data cells end up unstyled (style index 0) like the
# corresponding cells on every other quarter sheet.
$q4.Range("B2:G6").Style = "Normal"

# Apply the bold/centered/bordered header style (style index 2 in styles.xml)
# used by every other sheet's header + index column, by copying it from an
# existing sheet's header row.
$src = $wb.Worksheets.Item("2022-Q3")
$src.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$src.Range("A2").Copy()
$q4.Range("A2:A6").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert a new top data row for 2022-Q4
#    and push the existing quarters down by one row.
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Range("A2:D2").EntireRow.Insert()

# EntireRow.Insert() copies the format of the row above (the bold header),
# so strip it back to the unstyled default used by every other data row.
$totals.Range("A2:D2").Style = "Normal"

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 5
$totals.Range("D2").Value = 0.97

# The inserted row loses the index-column style; copy it back from the row
# beneath (still carrying the original "总计" index style).
$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3
$totals.Range("A6").Value = 4
$totals.Range("A7").Value = 5
$totals.Range("A8").Value = 6
$totals.Range("A9").Value = 7

# Keep "总计" as the active/selected sheet (matches the unchanged bookViews
# in the target - adding a sheet must not steal the active tab).
$totals.Activate()
